# Update DA plan ("Conditional indexation") column values with asset smoothing
# results on both the 15-year and 30-year pvERC sheets.

$wb = $excel.ActiveWorkbook

$ws15 = $wb.Worksheets.Item("pvERC_15y")
$ws15.Range("L2").Value = -13.613290668907252
$ws15.Range("L3").Value = -13.124547195962073
$ws15.Range("L4").Value = -11.567027137612262
$ws15.Range("L5").Value = -8.705518139597412
$ws15.Range("L6").Value = -1.336862033597086

$ws30 = $wb.Worksheets.Item("pvERC_30y")
$ws30.Range("L2").Value = -20.421180106816138
$ws30.Range("L3").Value = -18.002330808058638
$ws30.Range("L4").Value = -14.316716249182381
$ws30.Range("L5").Value = -12.626722370999898
$ws30.Range("L6").Value = -3.609838550652966
